# Auto-generated edit script: update FFXIV Adamantoise market price snapshot values
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (CUL is unchanged).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3363.5122
$ws.Range("I15").Value = 3363.5122
$ws.Range("K15").Value = 10090.5366
$ws.Range("M15").Value = -9921.536599999999
$ws.Range("H76").Value = 4850.5454
$ws.Range("I76").Value = 4958
$ws.Range("K76").Value = 4958
$ws.Range("M76").Value = -4643
$ws.Range("H79").Value = 4850.5454
$ws.Range("I79").Value = 4958
$ws.Range("K79").Value = 4958
$ws.Range("M79").Value = -3866
$ws.Range("H86").Value = 230773380
$ws.Range("J86").Value = 200002080
$ws.Range("L86").Value = 200002080
$ws.Range("N86").Value = -200004326
$ws.Range("H89").Value = 230773380
$ws.Range("J89").Value = 200002080
$ws.Range("L89").Value = 1000010400
$ws.Range("N89").Value = -1000021632
$ws.Range("H98").Value = 1188.2084
$ws.Range("I98").Value = 1184.381
$ws.Range("J98").Value = 1215
$ws.Range("K98").Value = 1184.381
$ws.Range("L98").Value = 1215
$ws.Range("M98").Value = 313.6189999999999
$ws.Range("N98").Value = -4211
$ws.Range("H100").Value = 3139.0588
$ws.Range("I100").Value = 1858
$ws.Range("J100").Value = 4277.778
$ws.Range("K100").Value = 1858
$ws.Range("L100").Value = 4277.778
$ws.Range("M100").Value = -1317
$ws.Range("N100").Value = -5359.778
$ws.Range("H106").Value = 3335442.8
$ws.Range("I106").Value = 3705714.2
$ws.Range("K106").Value = 3705714.2
$ws.Range("M106").Value = -3705083.2
$ws.Range("H107").Value = 3075.889
$ws.Range("I107").Value = 3373
$ws.Range("K107").Value = 3373
$ws.Range("M107").Value = -1453
$ws.Range("H122").Value = 1188.2084
$ws.Range("I122").Value = 1184.381
$ws.Range("J122").Value = 1215
$ws.Range("K122").Value = 3553.143
$ws.Range("L122").Value = 3645
$ws.Range("M122").Value = -1103.143
$ws.Range("N122").Value = -8545
$ws.Range("H131").Value = 3196.625
$ws.Range("I131").Value = 1862.3334
$ws.Range("K131").Value = 5587.0002
$ws.Range("M131").Value = -547.0002000000004
$ws.Range("H132").Value = 24307340
$ws.Range("I132").Value = 25091126
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 75273378
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -75270848
$ws.Range("N132").Value = -35057
$ws.Range("H137").Value = 3475525.8
$ws.Range("I137").Value = 3443.25
$ws.Range("K137").Value = 10329.75
$ws.Range("M137").Value = -7779.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = $null
$ws.Range("H45").Value = 4172.7354
$ws.Range("I45").Value = 3984.9333
$ws.Range("K45").Value = 3984.9333
$ws.Range("M45").Value = -3607.9333
$ws.Range("H102").Value = 1441.36
$ws.Range("J102").Value = 2639.8
$ws.Range("L102").Value = 2639.8
$ws.Range("N102").Value = -5883.8

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2779.3635
$ws.Range("I20").Value = 2208.111
$ws.Range("J20").Value = 5350
$ws.Range("K20").Value = 2208.111
$ws.Range("L20").Value = 5350
$ws.Range("M20").Value = -1961.111
$ws.Range("N20").Value = -5844
$ws.Range("H86").Value = 2845.9524
$ws.Range("I86").Value = 2665.9333
$ws.Range("J86").Value = 3296
$ws.Range("K86").Value = 2665.9333
$ws.Range("L86").Value = 3296
$ws.Range("M86").Value = -1542.9333
$ws.Range("N86").Value = -5542
$ws.Range("H89").Value = 2845.9524
$ws.Range("I89").Value = 2665.9333
$ws.Range("J89").Value = 3296
$ws.Range("K89").Value = 13329.6665
$ws.Range("L89").Value = 16480
$ws.Range("M89").Value = -7713.666500000001
$ws.Range("N89").Value = -27712

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 12500438
$ws.Range("I12").Value = 12500438
$ws.Range("K12").Value = 12500438
$ws.Range("M12").Value = -12500268
$ws.Range("H31").Value = 5753150.5
$ws.Range("I31").Value = 13890490
$ws.Range("J31").Value = 9146.706
$ws.Range("K31").Value = 13890490
$ws.Range("L31").Value = 9146.706
$ws.Range("M31").Value = -13890195
$ws.Range("N31").Value = -9736.706
$ws.Range("H34").Value = 5753150.5
$ws.Range("I34").Value = 13890490
$ws.Range("J34").Value = 9146.706
$ws.Range("K34").Value = 13890490
$ws.Range("L34").Value = 9146.706
$ws.Range("M34").Value = -13890288
$ws.Range("N34").Value = -9550.706
$ws.Range("H106").Value = 73553.336
$ws.Range("J106").Value = 73553.336
$ws.Range("L106").Value = 73553.336
$ws.Range("N106").Value = -76077.336
$ws.Range("H132").Value = 2576.5417
$ws.Range("I132").Value = 2253.524
$ws.Range("K132").Value = 6760.572
$ws.Range("M132").Value = -4230.572
$ws.Range("H134").Value = 83334670
$ws.Range("I134").Value = 83334670
$ws.Range("K134").Value = 250004010
$ws.Range("M134").Value = -250001475

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("H122").Value = 1600.6316
$ws.Range("I122").Value = 1363.5333
$ws.Range("K122").Value = 4090.5999
$ws.Range("M122").Value = -1640.5999
$ws.Range("H126").Value = 2750.1428
$ws.Range("I126").Value = 2097.5
$ws.Range("K126").Value = 6292.5
$ws.Range("M126").Value = -3822.5
$ws.Range("H136").Value = 61474.375
$ws.Range("J136").Value = 61474.375
$ws.Range("L136").Value = 184423.125
$ws.Range("N136").Value = -189523.125

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 109466.336
$ws.Range("J118").Value = 109466.336
$ws.Range("L118").Value = 109466.336
$ws.Range("N118").Value = -112780.336
$ws.Range("H132").Value = 4307.143
$ws.Range("J132").Value = 4880
$ws.Range("L132").Value = 14640
$ws.Range("N132").Value = -19700
$ws.Range("H136").Value = 14756.333
$ws.Range("I136").Value = 2332.5
$ws.Range("K136").Value = 6997.5
$ws.Range("M136").Value = -4447.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 29966.666
$ws.Range("J26").Value = 29950
$ws.Range("L26").Value = 29950
$ws.Range("N26").Value = -30536
$ws.Range("H54").Value = 30070
$ws.Range("I54").Value = 30070
$ws.Range("K54").Value = 30070
$ws.Range("M54").Value = -29550
$ws.Range("H74").Value = 35672.875
$ws.Range("J74").Value = 35672.875
$ws.Range("L74").Value = 35672.875
$ws.Range("N74").Value = -37544.875
$ws.Range("H77").Value = 35672.875
$ws.Range("J77").Value = 35672.875
$ws.Range("L77").Value = 107018.625
$ws.Range("N77").Value = -116378.625
$ws.Range("H81").Value = 4389.231
$ws.Range("I81").Value = 3666.75
$ws.Range("J81").Value = 4710.3335
$ws.Range("K81").Value = 7333.5
$ws.Range("L81").Value = 9420.666999999999
$ws.Range("M81").Value = -6272.5
$ws.Range("N81").Value = -11542.667
$ws.Range("H84").Value = 4389.231
$ws.Range("I84").Value = 3666.75
$ws.Range("J84").Value = 4710.3335
$ws.Range("K84").Value = 36667.5
$ws.Range("L84").Value = 47103.335
$ws.Range("M84").Value = -31363.5
$ws.Range("N84").Value = -57711.335
$ws.Range("H110").Value = 121000
$ws.Range("J110").Value = 121000
$ws.Range("L110").Value = 121000
$ws.Range("N110").Value = -129180
$ws.Range("H132").Value = 2661.9429
$ws.Range("I132").Value = 2294.96
$ws.Range("K132").Value = 6884.88
$ws.Range("M132").Value = -4354.88
$ws.Range("H133").Value = 44857.5
$ws.Range("J133").Value = 44857.5
$ws.Range("L133").Value = 44857.5
$ws.Range("N133").Value = -54977.5
$ws.Range("H136").Value = 2499.5217
$ws.Range("I136").Value = 2004.1111
$ws.Range("J136").Value = 4283
$ws.Range("K136").Value = 6012.3333
$ws.Range("L136").Value = 12849
$ws.Range("M136").Value = -3462.3333
$ws.Range("N136").Value = -17949
